# Actualización automática 2025-08-27 11:40:09
#
# Updates sales figures for advisor "HIDALGO HIDALGO PEDRO GUSTAVO" / client
# "FERRETERIAS FERRIGONZ SA" for the month of "agosto" (August) from 0 to
# 60.78 in the PORCELANATO product group, and propagates the resulting
# totals/percentages across the three report sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M10").Value = 60.78
$wsGrupo.Range("M23").Value = "4 de 21"

# --- Sheet "VENTA MENSUAL" --------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F10").Value = 60.78
$wsMensual.Range("F23").Value = 7532.75

# --- Sheet "CUMPLIMIENTO MENSUAL" ------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 6375.95
$wsCumplimiento.Range("E16").Value = 32400.52
$wsCumplimiento.Range("F16").Value = 0.1644283247030996

$wsCumplimiento.Range("D19").Value = 7532.75
$wsCumplimiento.Range("E19").Value = 51855.47762291769
$wsCumplimiento.Range("F19").Value = 0.1268391110748208
